# Apply cryptocurrency price/volume updates to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.365.16"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "2.928.66"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'597.58"
$ws.Range("E5").Value = "  +0.84%  "

$ws.Range("D6").Value = "'145.22"

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.96%  "

$ws.Range("D9").Value = "'7.02"
$ws.Range("E9").Value = "  +1.72%  "

$ws.Range("E10").Value = "  -2.48%  "

$ws.Range("E11").Value = "  -0.81%  "

$ws.Range("E12").Value = "  -1.29%  "

$ws.Range("D13").Value = "'33.50"
$ws.Range("E13").Value = "  -1.02%  "

$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("D15").Value = "3.415.29"
$ws.Range("E15").Value = "  -0.11%  "

$ws.Range("D16").Value = "61.347.60"
$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").Value = "2.930.38"
$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("D19").Value = "'432.41"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("E20").Value = "  +0.00%  "

$ws.Range("E21").Value = "  -1.23%  "

$ws.Range("E22").Value = "  -0.41%  "

$ws.Range("D23").Value = "'81.86"
$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("E24").Value = "  -2.04%  "

$ws.Range("D25").Value = "'2.18"
$ws.Range("E25").Value = "  -2.71%  "

$ws.Range("E26").Value = "  -2.29%  "

$ws.Range("E28").Value = "  -5.01%  "

$ws.Range("E29").Value = "  -0.47%  "

$ws.Range("D30").Value = "'6.91"
$ws.Range("E30").Value = "  -2.98%  "

$ws.Range("E31").Value = "  +1.13%  "

$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").Value = "0.0₃0885"
$ws.Range("E34").Value = "  +2.56%  "

$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("D36").Value = "'5.62"
$ws.Range("E36").Value = "  -0.34%  "

$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("E39").Value = "  -1.87%  "

$ws.Range("D40").Value = "'8.57"
$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("D41").Value = "'42.32"
$ws.Range("E41").Value = "  +6.26%  "

$ws.Range("E42").Value = "  -2.82%  "

$ws.Range("E43").Value = "  -0.43%  "

$ws.Range("D44").Value = "2.695.29"
$ws.Range("E44").Value = "  -0.86%  "

$ws.Range("D45").Value = "'133.84"
$ws.Range("E45").Value = "  +2.45%  "

$ws.Range("D46").Value = "'361.67"
$ws.Range("E46").Value = "  -4.15%  "

$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").Value = "'23.58"
$ws.Range("E48").Value = "  -3.00%  "

$ws.Range("E49").Value = "  -1.38%  "

$ws.Range("E50").Value = "  -2.17%  "

$ws.Range("E51").Value = "  -1.70%  "
